$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 82, shifting existing rows 82:111
# down to 83:112 (dimension grows from A1:R111 to A1:R112).
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row 82 with the new weekly price observation.
$ws.Cells.Item(82, 1).Value  = 1
$ws.Cells.Item(82, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(82, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(82, 4).Value  = 45016
$ws.Cells.Item(82, 5).Value  = 15
$ws.Cells.Item(82, 6).Value  = 100112040
$ws.Cells.Item(82, 7).Value  = "Cilantro"
$ws.Cells.Item(82, 8).Value  = "Sin especificar"
$ws.Cells.Item(82, 9).Value  = "Primera"
$ws.Cells.Item(82, 10).Value = 300
$ws.Cells.Item(82, 11).Value = 1800
$ws.Cells.Item(82, 12).Value = 2000
$ws.Cells.Item(82, 13).Value = 1900
$ws.Cells.Item(82, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(82, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(82, 16).Value = 950
$ws.Cells.Item(82, 17).Value = 2
$ws.Cells.Item(82, 18).Value = "Hortaliza"
